# "Sua CI chay TestNG, them test Bao mat, hieu nang"
# Update row 6 of the asset-issuance list with a new dispatch date,
# a new asset code, and the (existing) warehouse code "KTS-MOI".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "15/05/2025"
$ws.Range("C6").Value = "TS-007502"
$ws.Range("D6").Value = "KTS-MOI"
